$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Add a paragraph border with 5-twip spacing on every edge, and widen the
# left indent from 120 to 225 twips (6pt -> 11.25pt).
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Locate the ID placeholder text that lives in the first paragraph's first run.
$idRng = $p1.Range.Duplicate()
$idRng.Find.Execute("**ID__AFFARS_5322_topic_2__ID**", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

# Everything between the end of that text and the paragraph mark is the
# trailing space-only run; delete it outright (it is dropped in the target).
$paraMarkStart = $p1.Range.End - 1
$trailingRng = $d.Range($idRng.End, $paraMarkStart)
if ($trailingRng.Start -lt $trailingRng.End) {
    $trailingRng.Delete()
}

# Finally, rewrite the placeholder id itself.
$idRng.Text = "**ID__AFFARS_SUBPART_5322_1__ID**"
